# Applies the "Fixed Tora decision. GER country file WIP Added Werner von
# Blomberg as Field Marshal." edit to the Generals Skill levels workbook.
#
# The workbook has three sheets:
#   "Skill level" - the picker sheet with three dropdown-linked cells
#                    (A3 = skill level, A5 = personality 1, A7 = personality 2)
#   "Skills"      - lookup table for the skill-level dropdown
#   "Personality" - lookup table for the two personality dropdowns
#                    (left block B:F keyed by A, right block H:L keyed by G)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Skills sheet - direct data edit on the row for skill level 4
#    (A5 = 4): raise Attack/Defence/Logistics from 1 to 2, keep Planning at 1
# ---------------------------------------------------------------------
$skills = $wb.Worksheets.Item("Skills")
$skills.Range("C5").Value = 2
$skills.Range("D5").Value = 2
$skills.Range("F5").Value = 2
[void]$skills.Range("F5").Select()

# ---------------------------------------------------------------------
# 2. Personality sheet - direct data edits scattered through the two
#    lookup blocks (left B:F block and right H:L block)
# ---------------------------------------------------------------------
$personality = $wb.Worksheets.Item("Personality")

$personality.Range("D10").Value = 1
$personality.Range("F10").Value = 1

$personality.Range("D11").Value = 1
$personality.Range("J11").Value = 1
$personality.Range("L11").Value = 1

$personality.Range("C12").Value = 0
$personality.Range("D12").Value = 2
$personality.Range("E12").Value = 0
$personality.Range("J12").Value = 1
$personality.Range("L12").Value = 1

$personality.Range("I13").Value = 0
$personality.Range("J13").Value = 2
$personality.Range("K13").Value = 0

$personality.Range("F14").Value = -1

$personality.Range("L15").Value = -1

$personality.Range("E16").Value = 1

$personality.Range("K17").Value = 1

$personality.Range("C20").Value = 1

$personality.Range("D21").Value = 1
$personality.Range("F21").Value = 1
$personality.Range("I21").Value = 1

$personality.Range("D22").Value = 1
$personality.Range("E22").Value = 1
$personality.Range("J22").Value = 1
$personality.Range("L22").Value = 1

$personality.Range("J23").Value = 1
$personality.Range("K23").Value = 1

[void]$personality.Range("N25").Select()

# ---------------------------------------------------------------------
# 3. Skill level sheet - move the three dropdown-linked picker cells to
#    their new selections: skill level 3, personality-1 "8", personality-2 "22"
# ---------------------------------------------------------------------
$skillLevel = $wb.Worksheets.Item("Skill level")
[void]$skillLevel.Activate()

$skillLevel.Range("A3").Value = 3
$skillLevel.Range("A5").Value = 8
$skillLevel.Range("A7").Value = 22

[void]$skillLevel.Range("A5").Select()

[void]$wb.Application.Calculate()
